$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.269.46"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.550.65"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.14"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.771.64"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "1.550.05"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "28.263.20"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.36"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.67"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.30"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "0.0₃0672"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.79"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("E24").Value = "  -5.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.41"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.74"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("E31").Value = "  -4.50%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "1.382.10"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.770"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.39"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.59"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("E47").Value = "  -6.31%  "
$ws.Range("D48").Value = "1.685.71"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.28"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "41.56"
$ws.Range("E51").Value = "  +4.47%  "
